# GANCHO J.xlsx price list update: new date + updated prices for rows 29-37
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the list date shown in A1 (merged A1:D1)
$ws.Range("A1").Value = 45436

# Updated unit prices in column D
$ws.Range("D29").Value = 106.327
$ws.Range("D30").Value = 113.844
$ws.Range("D31").Value = 119.215
$ws.Range("D32").Value = 125.122
$ws.Range("D33").Value = 131.032
$ws.Range("D34").Value = 138.543
$ws.Range("D35").Value = 150.362
$ws.Range("D36").Value = 169.696
$ws.Range("D37").Value = 198.694
